$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("B9").Value = 6865285
$ws.Range("C9").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D9").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E9").Value = 45150.5
$ws.Range("F9").Value = "NK Igman Konjic"
$ws.Range("G9").Value = "Sloga"
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "H"
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 3.4
$ws.Range("M9").Value = 3.2
$ws.Range("N9").Value = 1.909
$ws.Range("O9").Value = 3.5
$ws.Range("P9").Value = 3.4
$ws.Range("Q9").Value = -0.5
$ws.Range("R9").Value = 1.95
$ws.Range("S9").Value = 1.85
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 1.85
$ws.Range("V9").Value = 1.95
$ws.Range("W9").Value = 0.909
$ws.Range("X9").Value = -1
$ws.Range("Y9").Value = -1
$ws.Range("Z9").Value = 0.95
$ws.Range("AA9").Value = -1
$ws.Range("AB9").Value = -1
$ws.Range("AC9").Value = 0.95

# Row 10
$ws.Range("B10").Value = 6865281
$ws.Range("C10").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D10").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E10").Value = 45150.5
$ws.Range("F10").Value = "GOSK Gabela"
$ws.Range("G10").Value = "Zvijezda 09"
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "H"
$ws.Range("K10").Value = 1.75
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 3.5
$ws.Range("N10").Value = 1.75
$ws.Range("O10").Value = 4
$ws.Range("P10").Value = 3.4
$ws.Range("Q10").Value = -0.5
$ws.Range("R10").Value = 1.8
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 2.5
$ws.Range("U10").Value = 1.85
$ws.Range("V10").Value = 1.95
$ws.Range("W10").Value = 0.75
$ws.Range("X10").Value = -1
$ws.Range("Y10").Value = -1
$ws.Range("Z10").Value = 0.8
$ws.Range("AA10").Value = -1
$ws.Range("AB10").Value = -1
$ws.Range("AC10").Value = 0.95

# Row 36
$ws.Range("B36").Value = 6865299
$ws.Range("C36").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D36").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E36").Value = 45186.61458333334
$ws.Range("F36").Value = "Siroki Brijeg"
$ws.Range("G36").Value = "Zvijezda 09"
$ws.Range("H36").Value = 2
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = "H"
$ws.Range("K36").Value = 1.25
$ws.Range("L36").Value = 5.5
$ws.Range("M36").Value = 8
$ws.Range("N36").Value = 1.4
$ws.Range("O36").Value = 4.75
$ws.Range("P36").Value = 5.75
$ws.Range("Q36").Value = -1.25
$ws.Range("R36").Value = 1.9
$ws.Range("S36").Value = 1.9
$ws.Range("T36").Value = 2.75
$ws.Range("U36").Value = 1.85
$ws.Range("V36").Value = 1.95
$ws.Range("W36").Value = 0.3999999999999999
$ws.Range("X36").Value = -1
$ws.Range("Y36").Value = -1
$ws.Range("Z36").Value = -0.5
$ws.Range("AA36").Value = 0.45
$ws.Range("AB36").Value = 0.425
$ws.Range("AC36").Value = -0.5

# Row 37
$ws.Range("B37").Value = 6864629
$ws.Range("C37").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D37").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E37").Value = 45186.61458333334
$ws.Range("F37").Value = "Borac Banja Luka"
$ws.Range("G37").Value = "NK Posusje"
$ws.Range("H37").Value = 1
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = "H"
$ws.Range("K37").Value = 1.363
$ws.Range("L37").Value = 4.5
$ws.Range("M37").Value = 6.5
$ws.Range("N37").Value = 1.363
$ws.Range("O37").Value = 4.2
$ws.Range("P37").Value = 6.5
$ws.Range("Q37").Value = -1.25
$ws.Range("R37").Value = 1.95
$ws.Range("S37").Value = 1.85
$ws.Range("T37").Value = 2.5
$ws.Range("U37").Value = 1.925
$ws.Range("V37").Value = 1.875
$ws.Range("W37").Value = 0.363
$ws.Range("X37").Value = -1
$ws.Range("Y37").Value = -1
$ws.Range("Z37").Value = -0.5
$ws.Range("AA37").Value = 0.425
$ws.Range("AB37").Value = -1
$ws.Range("AC37").Value = 0.875

# Row 76
$ws.Range("B76").Value = 6865328
$ws.Range("C76").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D76").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E76").Value = 45235.375
$ws.Range("F76").Value = "Siroki Brijeg"
$ws.Range("G76").Value = "NK Posusje"
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = "D"
$ws.Range("K76").Value = 2
$ws.Range("L76").Value = 3
$ws.Range("M76").Value = 3.5
$ws.Range("N76").Value = 2.1
$ws.Range("O76").Value = 3
$ws.Range("P76").Value = 3.3
$ws.Range("Q76").Value = -0.25
$ws.Range("R76").Value = 1.825
$ws.Range("S76").Value = 1.975
$ws.Range("T76").Value = 2
$ws.Range("U76").Value = 1.825
$ws.Range("V76").Value = 1.975
$ws.Range("W76").Value = -1
$ws.Range("X76").Value = 2
$ws.Range("Y76").Value = -1
$ws.Range("Z76").Value = -0.5
$ws.Range("AA76").Value = 0.4875
$ws.Range("AB76").Value = 0
$ws.Range("AC76").Value = -0.0

# Row 77
$ws.Range("B77").Value = 6865377
$ws.Range("C77").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D77").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E77").Value = 45235.375
$ws.Range("F77").Value = "Zrinjski Mostar"
$ws.Range("G77").Value = "FK Tuzla City"
$ws.Range("H77").Value = 3
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = "H"
$ws.Range("K77").Value = 1.333
$ws.Range("L77").Value = 5
$ws.Range("M77").Value = 6
$ws.Range("N77").Value = 1.166
$ws.Range("O77").Value = 6.5
$ws.Range("P77").Value = 13
$ws.Range("Q77").Value = -2
$ws.Range("R77").Value = 1.9
$ws.Range("S77").Value = 1.9
$ws.Range("T77").Value = 3.25
$ws.Range("U77").Value = 1.95
$ws.Range("V77").Value = 1.85
$ws.Range("W77").Value = 0.1659999999999999
$ws.Range("X77").Value = -1
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = 0
$ws.Range("AA77").Value = -0.0
$ws.Range("AB77").Value = 0.95
$ws.Range("AC77").Value = -1

# Row 99
$ws.Range("B99").Value = 6864639
$ws.Range("C99").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D99").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E99").Value = 45269.375
$ws.Range("F99").Value = "Zvijezda 09"
$ws.Range("G99").Value = "Borac Banja Luka"
$ws.Range("H99").Value = 1
$ws.Range("I99").Value = 2
$ws.Range("J99").Value = "A"
$ws.Range("K99").Value = 11
$ws.Range("L99").Value = 6
$ws.Range("M99").Value = 1.2
$ws.Range("N99").Value = 10
$ws.Range("O99").Value = 6.5
$ws.Range("P99").Value = 1.181
$ws.Range("Q99").Value = 2
$ws.Range("R99").Value = 1.825
$ws.Range("S99").Value = 1.975
$ws.Range("T99").Value = 3
$ws.Range("U99").Value = 1.9
$ws.Range("V99").Value = 1.9
$ws.Range("W99").Value = -1
$ws.Range("X99").Value = -1
$ws.Range("Y99").Value = 0.181
$ws.Range("Z99").Value = 0.825
$ws.Range("AA99").Value = -1
$ws.Range("AB99").Value = 0
$ws.Range("AC99").Value = -0.0

# Row 100
$ws.Range("B100").Value = 6865343
$ws.Range("C100").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D100").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E100").Value = 45269.375
$ws.Range("F100").Value = "Sloga"
$ws.Range("G100").Value = "NK Posusje"
$ws.Range("H100").Value = 1
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = "H"
$ws.Range("K100").Value = 1.909
$ws.Range("L100").Value = 3.3
$ws.Range("M100").Value = 3.5
$ws.Range("N100").Value = 2.2
$ws.Range("O100").Value = 2.8
$ws.Range("P100").Value = 3.3
$ws.Range("Q100").Value = -0.25
$ws.Range("R100").Value = 1.95
$ws.Range("S100").Value = 1.85
$ws.Range("T100").Value = 1.75
$ws.Range("U100").Value = 1.875
$ws.Range("V100").Value = 1.925
$ws.Range("W100").Value = 1.2
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = 0.95
$ws.Range("AA100").Value = -1
$ws.Range("AB100").Value = -1
$ws.Range("AC100").Value = 0.925

# Row 122
$ws.Range("B122").Value = 6865363
$ws.Range("C122").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D122").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E122").Value = 45353.375
$ws.Range("F122").Value = "NK Igman Konjic"
$ws.Range("G122").Value = "Siroki Brijeg"
$ws.Range("H122").Value = 1
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = "H"
$ws.Range("K122").Value = 2
$ws.Range("L122").Value = 3.3
$ws.Range("M122").Value = 3.25
$ws.Range("N122").Value = 2.3
$ws.Range("O122").Value = 3.2
$ws.Range("P122").Value = 2.75
$ws.Range("Q122").Value = -0.25
$ws.Range("R122").Value = 2.05
$ws.Range("S122").Value = 1.75
$ws.Range("T122").Value = 2
$ws.Range("U122").Value = 1.9
$ws.Range("V122").Value = 1.9
$ws.Range("W122").Value = 1.3
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = 1.05
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = -1
$ws.Range("AC122").Value = 0.8999999999999999

# Row 123
$ws.Range("B123").Value = 6865381
$ws.Range("C123").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D123").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E123").Value = 45353.375
$ws.Range("F123").Value = "FK Tuzla City"
$ws.Range("G123").Value = "Zvijezda 09"
$ws.Range("H123").Value = 2
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = "H"
$ws.Range("K123").Value = 1.666
$ws.Range("L123").Value = 3.6
$ws.Range("M123").Value = 4.333
$ws.Range("N123").Value = 1.5
$ws.Range("O123").Value = 4
$ws.Range("P123").Value = 5.25
$ws.Range("Q123").Value = -1
$ws.Range("R123").Value = 1.925
$ws.Range("S123").Value = 1.875
$ws.Range("T123").Value = 2.5
$ws.Range("U123").Value = 1.8
$ws.Range("V123").Value = 2
$ws.Range("W123").Value = 0.5
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = -1
$ws.Range("Z123").Value = 0.925
$ws.Range("AA123").Value = -1
$ws.Range("AB123").Value = -1
$ws.Range("AC123").Value = 1

# Row 124
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = "D"
$ws.Range("K124").Value = 5.75
$ws.Range("L124").Value = 4.5
$ws.Range("M124").Value = 1.363
$ws.Range("N124").Value = 5
$ws.Range("O124").Value = 3.8
$ws.Range("P124").Value = 1.533
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = 1.8
$ws.Range("S124").Value = 2
$ws.Range("T124").Value = 2
$ws.Range("U124").Value = 1.75
$ws.Range("V124").Value = 2.05
$ws.Range("W124").Value = -1
$ws.Range("X124").Value = 2.8
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = 0.8
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = 0
$ws.Range("AC124").Value = -0.0

# Row 125
$ws.Range("A125").Value = 123
$ws.Range("B125").Value = 6865362
$ws.Range("C125").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D125").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E125").Value = 45354.41666666666
$ws.Range("F125").Value = "NK Posusje"
$ws.Range("G125").Value = "Velez Mostar"
$ws.Range("K125").Value = 3.4
$ws.Range("L125").Value = 2.9
$ws.Range("M125").Value = 2.15
$ws.Range("N125").Value = 2.625
$ws.Range("O125").Value = 2.8
$ws.Range("P125").Value = 2.625
$ws.Range("Q125").Value = 0
$ws.Range("R125").Value = 1.9
$ws.Range("S125").Value = 1.9
$ws.Range("T125").Value = 1.75
$ws.Range("U125").Value = 1.8
$ws.Range("V125").Value = 2
$ws.Range("W125").Value = 0
$ws.Range("X125").Value = 0
$ws.Range("Y125").Value = 0
$ws.Range("Z125").Value = 0
$ws.Range("AA125").Value = 0

# Row 126
$ws.Range("A126").Value = 124
$ws.Range("B126").Value = 6865364
$ws.Range("C126").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D126").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E126").Value = 45354.61458333334
$ws.Range("F126").Value = "Zeljeznicar"
$ws.Range("G126").Value = "FK Sarajevo"
$ws.Range("K126").Value = 3.25
$ws.Range("L126").Value = 3
$ws.Range("M126").Value = 2.2
$ws.Range("N126").Value = 3.1
$ws.Range("O126").Value = 3
$ws.Range("P126").Value = 2.3
$ws.Range("Q126").Value = 0.25
$ws.Range("R126").Value = 1.8
$ws.Range("S126").Value = 2
$ws.Range("T126").Value = 2
$ws.Range("U126").Value = 1.85
$ws.Range("V126").Value = 1.95
$ws.Range("W126").Value = 0
$ws.Range("X126").Value = 0
$ws.Range("Y126").Value = 0
$ws.Range("Z126").Value = 0
$ws.Range("AA126").Value = 0

# Copy formatting for new rows 125 and 126 (column A bold/border/center style, column E date format)
$ws.Range("A124").Copy()
$ws.Range("A125:A126").PasteSpecial(-4122)
$ws.Range("E124").Copy()
$ws.Range("E125:E126").PasteSpecial(-4122)
$excel.CutCopyMode = 0
